$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the explicit "General" number format to B1 (adds a new cellXfs entry
# with applyNumberFormat="1", matching the fix that normalizes B1's format).
$ws.Range("B1").NumberFormat = "General"

# Updated computed values (bugfix in time computing / Greedy algorithm).
$ws.Range("AE1").Value = 9.8584060415702996
$ws.Range("AI1").Value = 9.8621617420716507
$ws.Range("AL1").Value = 5.3002636771839402
$ws.Range("AP1").Value = 5.3002636771839402
$ws.Range("B1").Value = 1003.59505890593
$ws.Range("C1").Value = 9.8443438331496207
$ws.Range("G1").Value = 9.8418471971860697
$ws.Range("Q1").Value = 9.8584060415702996
$ws.Range("T1").Value = 0.168095472122096
$ws.Range("U1").Value = 9.9381617420716495
$ws.Range("AE2").Value = 20.246913780813198
$ws.Range("AI2").Value = 20.142668575374898
$ws.Range("AL2").Value = 9.9033685267321747
$ws.Range("AP2").Value = 9.9033685267321747
$ws.Range("B2").Value = 2007.19
$ws.Range("C2").Value = 19.762205894692976
$ws.Range("G2").Value = 19.783816959691134
$ws.Range("Q2").Value = 20.246913780813198
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0.30621881977772297
$ws.Range("U2").Value = 20.243668575374901
$ws.Range("V2").Value = 0
$ws.Range("AE3").Value = 29.971552574152799
$ws.Range("AI3").Value = 29.830050650016101
$ws.Range("AL3").Value = 12.407563508323898
$ws.Range("AP3").Value = 12.407563508323898
$ws.Range("B3").Value = 3010.7849999999999
$ws.Range("C3").Value = 29.115972251621969
$ws.Range("G3").Value = 29.145853136844799
$ws.Range("Q3").Value = 33.6820803978732
$ws.Range("R3").Value = 3.6005505040752901
$ws.Range("S3").Value = 1
$ws.Range("T3").Value = 0.25953717764963702
$ws.Range("U3").Value = 33.5153119842563
$ws.Range("V3").Value = 1
$ws.Range("AE4").Value = 40.579156675758099
$ws.Range("AI4").Value = 40.389720018337499
$ws.Range("AL4").Value = 16.860545565823742
$ws.Range("AP4").Value = 16.860545565823742
$ws.Range("B4").Value = 4014.3802356237202
$ws.Range("C4").Value = 41.924944702572205
$ws.Range("G4").Value = 41.918769534202987
$ws.Range("Q4").Value = 44.5139527215917
$ws.Range("R4").Value = 3.60279407845462
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 0.55463936423556404
$ws.Range("U4").Value = 44.492482389639697
$ws.Range("V4").Value = 1
$ws.Range("AE5").Value = 51.806655171450416
$ws.Range("AI5").Value = 51.63602472628682
$ws.Range("AL5").Value = 20.319508511507696
$ws.Range("AP5").Value = 20.319508511507696
$ws.Range("B5").Value = 5017.9750000000004
$ws.Range("C5").Value = 52.814924323620332
$ws.Range("G5").Value = 52.763487806438086
$ws.Range("Q5").Value = 59.013022503787198
$ws.Range("R5").Value = 7.2063673323367903
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 0.47483142729368699
$ws.Range("U5").Value = 58.896324587335798
$ws.Range("V5").Value = 2
$ws.Range("AE6").Value = 63.58193729589297
$ws.Range("AI6").Value = 63.309559718920667
$ws.Range("AL6").Value = 23.661184528677854
$ws.Range("AP6").Value = 23.661184528677854
$ws.Range("B6").Value = 6021.57
$ws.Range("C6").Value = 66.441826425685477
$ws.Range("G6").Value = 65.722314129195553
$ws.Range("Q6").Value = 72.739582554584601
$ws.Range("R6").Value = 10.8113858788711
$ws.Range("S6").Value = 3
$ws.Range("T6").Value = 0.81762305768054699
$ws.Range("U6").Value = 72.736482035728102
$ws.Range("V6").Value = 3
$ws.Range("AE7").Value = 71.275137071940804
$ws.Range("AI7").Value = 71.08841213164861
$ws.Range("AL7").Value = 27.182327681147452
$ws.Range("AP7").Value = 27.182327681147452
$ws.Range("B7").Value = 7025.165
$ws.Range("C7").Value = 76.751976545525295
$ws.Range("G7").Value = 75.166709965261504
$ws.Range("Q7").Value = 86.239626855793304
$ws.Range("R7").Value = 14.429689694999899
$ws.Range("S7").Value = 4
$ws.Range("T7").Value = 0.67079966751631204
$ws.Range("U7").Value = 85.936981971805906
$ws.Range("V7").Value = 4
$ws.Range("AE8").Value = 83.233399625618802
$ws.Range("AI8").Value = 82.751206203536711
$ws.Range("AL8").Value = 34.878759329679497
$ws.Range("AP8").Value = 34.878759329679497
$ws.Range("B8").Value = 8028.76
$ws.Range("C8").Value = 90.488391899777923
$ws.Range("G8").Value = 89.579458319030508
$ws.Range("Q8").Value = 97.670613402044097
$ws.Range("R8").Value = 14.437213776425301
$ws.Range("S8").Value = 4
$ws.Range("T8").Value = 0.76349608810185599
$ws.Range("U8").Value = 97.605318861791702
$ws.Range("V8").Value = 4
$ws.Range("AE9").Value = 95.656717140325796
$ws.Range("AI9").Value = 96.109047046121802
$ws.Range("AL9").Value = 42.033180460678423
$ws.Range("AP9").Value = 42.033180460678423
$ws.Range("B9").Value = 9032.3549999999996
$ws.Range("C9").Value = 97.767615081483342
$ws.Range("G9").Value = 96.015942324297782
$ws.Range("Q9").Value = 124.115303890785
$ws.Range("R9").Value = 18.064179306943402
$ws.Range("S9").Value = 5
$ws.Range("T9").Value = 0.77696965770298698
$ws.Range("U9").Value = 123.67664374211
$ws.Range("V9").Value = 5
$ws.Range("AE10").Value = 172.62497268322946
$ws.Range("AI10").Value = 172.75331597934044
$ws.Range("AL10").Value = 48.945727191344986
$ws.Range("AP10").Value = 48.945727191344986
$ws.Range("B10").Value = 10035.950000000001
$ws.Range("C10").Value = 109.5009830679254
$ws.Range("G10").Value = 108.51546685769162
$ws.Range("Q10").Value = 391.97301839833398
$ws.Range("R10").Value = 21.6745402885741
$ws.Range("S10").Value = 6
$ws.Range("T10").Value = 0.691298614282576
$ws.Range("U10").Value = 392.11048864042999
$ws.Range("V10").Value = 6

# Page setup: orientation now explicitly set to portrait.
$ws.PageSetup.Orientation = 1

# Selection moved to B14.
$ws.Range("B14").Select()
